# Fixed workflow: drop the first 4 data rows (Cutoff 0-3) from each results
# sheet and renumber the remaining "Cutoff" column starting again at 0.
# Net effect on disk: rows shift up by 4 and the sheet's used range shrinks
# from A1:C20 to A1:C16.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Delete the 4 obsolete data rows (rows 2-5); Excel shifts rows 6-20 up.
    $ws.Range("A2:A5").EntireRow.Delete() | Out-Null

    # After the delete, data that used to be in rows 6-20 now lives in
    # rows 2-16. Renumber column A (Cutoff) back to a 0-based sequence.
    for ($i = 0; $i -le 14; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
